# Auto-generated edit script applying crypto price/volume/coin-order updates
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.001.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.82%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.818.53"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.10%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9983"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.36%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.81%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9970"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.46%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4714"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.08%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3713"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.98%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07382"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.16%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8743"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.77%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.48"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.52%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.878.42"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.40%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.368"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.42%  "

# Row 14
$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.36%  "

# Row 15
$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07079"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.16%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.523"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.96%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9997"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.27%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008722"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.28%  "

# Row 19
$ws.Range("E19").Value = "  +0.05%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.04%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.013.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.79%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.348"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.47%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.58"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.98%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.101.62"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.66%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.894"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.06%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.95"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.08%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.46"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.52%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.161"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.95%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.309"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.57%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.64"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.10%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08956"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.80%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7631"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.63%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.164"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.24%  "

# Row 34
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.483"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.94%  "

# Row 35
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.925"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.08%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9984"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.28%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.097"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.39%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01958"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.19%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05269"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.50%  "

# Row 40
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5380"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.71%  "

# Row 41
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.927"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.54%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.238"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.47%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.383"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.19%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1665"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.05%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.490"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.22%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4971"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.50%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.36"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.00%  "

# Row 48
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.681"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.65%  "

# Row 49
$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.9979"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.34%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "103.39"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.67%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06304"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.36%  "
